$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 18:20"

# Update country statistics that changed between the two data snapshots
# (values below are the authoritative post-edit numbers per column B..H)
# Row 4: Estados Unidos
$ws.Range("B4").Value = 174750
$ws.Range("C4").Value = 10962
$ws.Range("D4").Value = 6215
$ws.Range("E4").Value = 165133
$ws.Range("F4").Value = 3893
$ws.Range("G4").Value = 261
$ws.Range("H4").Value = 3402

# Row 5: Italia
$ws.Range("B5").Value = 105792
$ws.Range("C5").Value = 4053
$ws.Range("D5").Value = 15729
$ws.Range("E5").Value = 77635
$ws.Range("F5").Value = 4023
$ws.Range("G5").Value = 837
$ws.Range("H5").Value = 12428

# Row 16: Austria
$ws.Range("B16").Value = 10088
$ws.Range("C16").Value = 470
$ws.Range("D16").Value = 1095
$ws.Range("E16").Value = 8865
$ws.Range("F16").Value = 198
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 128

# Row 33: Filipinas
$ws.Range("B33").Value = 2178
$ws.Range("C33").Value = 190
$ws.Range("D33").Value = 80
$ws.Range("E33").Value = 2075
$ws.Range("F33").Value = 31
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 23

# Row 34: Luxemburgo
$ws.Range("B34").Value = 2084
$ws.Range("C34").Value = 538
$ws.Range("D34").Value = 49
$ws.Range("E34").Value = 1947
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 88

# Row 57: Hong Kong
$ws.Range("B57").Value = 716
$ws.Range("C57").Value = 132
$ws.Range("D57").Value = 46
$ws.Range("E57").Value = 626
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 44

# Row 58: Crucero
$ws.Range("B58").Value = 714
$ws.Range("C58").Value = 31
$ws.Range("D58").Value = 128
$ws.Range("E58").Value = 582
$ws.Range("F58").Value = 5
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 4

# Row 59: Catar
$ws.Range("B59").Value = 712
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 603
$ws.Range("E59").Value = 99
$ws.Range("F59").Value = 15
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 10

# Row 60: Emiratos Arabes Unidos
$ws.Range("B60").Value = 693
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 51
$ws.Range("E60").Value = 641
$ws.Range("F60").Value = 6
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1

# Row 61: Egipto
$ws.Range("B61").Value = 664
$ws.Range("C61").Value = 53
$ws.Range("D61").Value = 61
$ws.Range("E61").Value = 597
$ws.Range("F61").Value = 2
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 6

# Row 62: Nueva Zelanda
$ws.Range("B62").Value = 656
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 150
$ws.Range("E62").Value = 465
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 41

# Row 63: Irak
$ws.Range("B63").Value = 647
$ws.Range("C63").Value = 58
$ws.Range("D63").Value = 74
$ws.Range("E63").Value = 572
$ws.Range("F63").Value = 2
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 1

# Row 64: Argelia
$ws.Range("B64").Value = 630
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 152
$ws.Range("E64").Value = 432
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 46

# Row 105: Sri Lanka
$ws.Range("B105").Value = 142
$ws.Range("C105").Value = 20
$ws.Range("D105").Value = 17
$ws.Range("E105").Value = 123
$ws.Range("F105").Value = 5
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 2

# Row 123: Liechtenstein
$ws.Range("B123").Value = 68
$ws.Range("C123").Value = 6
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 68
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0

# Row 156: Haiti
$ws.Range("B156").Value = 15
$ws.Range("C156").Value = 1
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 14
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 1

# Row 157: Birmania
$ws.Range("B157").Value = 15
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 1
$ws.Range("E157").Value = 14
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

# Row 172: Santa Lucia
$ws.Range("B172").Value = 9
$ws.Range("C172").Value = 3
$ws.Range("D172").Value = 1
$ws.Range("E172").Value = 8
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

# Row 173: Libia
$ws.Range("B173").Value = 9
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 8
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

# Row 178: Zimbabue
$ws.Range("B178").Value = 8
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179: Guyana
$ws.Range("B179").Value = 8
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

# Row 180: Antigua y Barbuda
$ws.Range("B180").Value = 8
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 7
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 1

# Row 182: Sudan
$ws.Range("B182").Value = 7
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 7
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0

# Row 184: Santa Sede
$ws.Range("B184").Value = 7
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 1
$ws.Range("E184").Value = 4
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 2

# Row 186: Benin
$ws.Range("B186").Value = 6
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 6
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187: San Bartolome
$ws.Range("B187").Value = 6
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 5
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 1

# Row 188: Cabo Verde
$ws.Range("B188").Value = 6
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 5
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 205: Sierra Leona
$ws.Range("B205").Value = 1
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 1
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Row 206: Papua Nueva Guinea
$ws.Range("B206").Value = 1
$ws.Range("C206").Value = 1
$ws.Range("D206").Value = 0
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

